$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Panel de información" + [_GoBack bookmark] + " de lotes"
#          -> single run "Panel de información de lotes" (bookmark removed)
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "Panel de información de lotes", $true, $false, $false, $false, $false,
    $true, 1, $false, "Panel de información de lotes", 2)

# ---------------------------------------------------------------------------
# Edit 2: the second ": NO." (the one that follows "Extends") becomes a long
#          explanation ending in " ." with a fresh _GoBack bookmark wrapping
#          the single space right before the final period.
# ---------------------------------------------------------------------------
$locate = $d.Content.Duplicate
$locate.Find.Execute("Extends: NO.")

$target = $d.Range($locate.End - 5, $locate.End)
$target.Text = ": Abrir lote  RF [17.6], Cerrar lote  RF [17.5]"
$target.Collapse(0)

$bmStart = $target.End
$target.InsertAfter(" .")
$bmEnd = $target.End - 1

$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Edit 3: "Crear lote" -> "Panel de información de lotes"
# ---------------------------------------------------------------------------
$found3 = $d.Content.Find.Execute(
    "Crear lote", $true, $false, $false, $false, $false,
    $true, 1, $false, "Panel de información de lotes", 2)
